$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "134.25") are not silently converted to numbers, matching the
# workbook's original inline-string representation.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = '61.330.93'
$ws.Range("D3").Value = '3.375.36'
$ws.Range("D5").Value = '405.73'
$ws.Range("D6").Value = '134.25'
$ws.Range("D7").Value = '0.592'
$ws.Range("D9").Value = '0.670'
$ws.Range("D10").Value = '0.120'
$ws.Range("D11").Value = '42.71'
$ws.Range("D13").Value = '3.894.57'
$ws.Range("D15").Value = '19.72'
$ws.Range("D16").Value = '3.356.25'
$ws.Range("D17").Value = '61.334.13'
$ws.Range("D21").Value = '3.21'
$ws.Range("D22").Value = '85.06'
$ws.Range("D23").Value = '314.69'
$ws.Range("D24").Value = '12.82'
$ws.Range("D25").Value = '3.13'
$ws.Range("D26").Value = '4.78'
$ws.Range("D27").Value = '8.31'
$ws.Range("D28").Value = '29.49'
$ws.Range("D29").Value = '7.61'
$ws.Range("D31").Value = '2.66'
$ws.Range("D32").Value = '0.170'
$ws.Range("D33").Value = '11.37'
$ws.Range("D35").Value = '41.06'
$ws.Range("D36").Value = '0.0481'
$ws.Range("D37").Value = '51.85'
$ws.Range("D38").Value = '0.998'
$ws.Range("D41").Value = '139.83'
$ws.Range("D44").Value = '0.296'
$ws.Range("D45").Value = '4.02'
$ws.Range("D46").Value = '16.70'
$ws.Range("D48").Value = '21.40'
$ws.Range("D49").Value = '2.117.18'
$ws.Range("D51").Value = '1.93'

# --- Coin / Link / Volume(1h) updates ---
$ws.Range("E2").Value = '  -2.37%  '
$ws.Range("E3").Value = '  -2.73%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("E6").Value = '  +7.88%  '
$ws.Range("E7").Value = '  -0.75%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  -5.74%  '
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("E13").Value = '  -3.05%  '
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("E20").Value = '  -5.72%  '
$ws.Range("E21").Value = '  -4.15%  '
$ws.Range("E22").Value = '  +2.58%  '
$ws.Range("E23").Value = '  -1.62%  '
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("E26").Value = '  +11.00%  '
$ws.Range("E27").Value = '  +5.39%  '
$ws.Range("E28").Value = '  -5.33%  '
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E31").Value = '  +3.78%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E32").Value = '  -1.72%  '
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -2.90%  '
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("E40").Value = '  -3.62%  '
$ws.Range("E41").Value = '  +3.77%  '
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("E44").Value = '  +3.87%  '
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("E48").Value = '  -3.56%  '
$ws.Range("E49").Value = '  -4.07%  '
$ws.Range("E50").Value = '  -5.44%  '
$ws.Range("E51").Value = '  +0.69%  '
